# Insert a new data row at row 226 (pushing the existing rows 226:304 down
# to 227:305, and extending the sheet's used range to T305), then populate
# the newly inserted row with the new "Granada" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 226..304 down by one row.
$ws.Rows.Item(226).Insert()

# Fill in the new row 226 with the new record's data.
$ws.Cells.Item(226, 1).Value  = 10
$ws.Cells.Item(226, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(226, 3).Value  = "La Araucanía"
$ws.Cells.Item(226, 4).Value  = 45215
$ws.Cells.Item(226, 5).Value  = 9
$ws.Cells.Item(226, 6).Value  = "Fruta"
$ws.Cells.Item(226, 7).Value  = 100104
$ws.Cells.Item(226, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(226, 9).Value  = 100104001
$ws.Cells.Item(226, 10).Value = "Granada"
$ws.Cells.Item(226, 11).Value = "Wonderfull"
$ws.Cells.Item(226, 12).Value = "Primera"
$ws.Cells.Item(226, 13).Value = 110
$ws.Cells.Item(226, 14).Value = 17000
$ws.Cells.Item(226, 15).Value = 17000
$ws.Cells.Item(226, 16).Value = 17000
$ws.Cells.Item(226, 17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(226, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(226, 19).Value = 1700
$ws.Cells.Item(226, 20).Value = 10
